$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Insert two blank rows at 33:34. This pushes the old rows 31-43
#    down to 33-45 (old row 31 -> 33, old row 32 -> 34, old row 36
#    "RAZEM" header -> 38, and the hydrogen table old rows 39-43
#    -> 41-45), which already matches the target layout for
#    everything from row 35 onward without any further edits.
# ------------------------------------------------------------------
$ws.Rows("33:34").Insert()

# ------------------------------------------------------------------
# 2) The two new blank rows (33:34) should become exact copies of
#    what is still sitting in rows 31 and 32 (the old, untouched
#    content), so duplicate that content+formatting down first.
# ------------------------------------------------------------------
$ws.Range("B31:F31").Copy()
$ws.Range("B33:F33").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B33:F33").PasteSpecial(-4163)   # xlPasteValues

$ws.Range("B32:F32").Copy()
$ws.Range("B34:F34").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B34:F34").PasteSpecial(-4163)   # xlPasteValues

$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 3) Row 30 switches from the "white" row style to the "grey" row
#    style used by rows 29/31/33 before its values are changed.
# ------------------------------------------------------------------
$ws.Range("B29:D29").Copy()
$ws.Range("B30:D30").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Row 32's E cell picks up the regular grey-row numeric style too.
$ws.Range("E31").Copy()
$ws.Range("E32").PasteSpecial(-4122)       # xlPasteFormats
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# 4) Now overwrite rows 28-32 with the new NCAP/CAP bound schedule
#    for ELE_NEW_NUC_PWR.
# ------------------------------------------------------------------

# Row 28: FX/NCAP_BND -> LO/CAP_BND (2030, 3 unchanged)
$ws.Range("B28").Value = "LO"
$ws.Range("C28").Value = "CAP_BND"

# Row 29: UP -> LO (2035 unchanged), value 10 -> 1.5
$ws.Range("B29").Value = "LO"
$ws.Range("E29").Value = 1.5

# Row 30: becomes UP/NCAP_BND, 2035, 10
$ws.Range("B30").Value = "UP"
$ws.Range("C30").Value = "NCAP_BND"
$ws.Range("D30").Value = 2035
$ws.Range("E30").Value = 10

# Row 31: LO/CAP_BND stays, year 2050 -> 2040, value 15 -> 5
$ws.Range("D31").Value = 2040
$ws.Range("E31").Value = 5

# Row 32: CAP_BND -> NCAP_BND, year 2050 -> 2040, value 25 -> 5
$ws.Range("C32").Value = "NCAP_BND"
$ws.Range("D32").Value = 2040
$ws.Range("E32").Value = 5

# ------------------------------------------------------------------
# 5) Sheet view bookkeeping to match the saved state.
# ------------------------------------------------------------------
$ws.Range("M29").Select()
$ws.Application.ActiveWindow.ScrollRow = 11
